# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output data (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Row -> new value for column F, shared by both affected sheets.
$updates = @{
    2  = 185
    3  = 432
    4  = 12560
    5  = 1280
    6  = 149
    7  = 33
    9  = 162
    10 = 203
    11 = 455
    17 = 4234
    18 = 98
    19 = 17
    20 = 944
    21 = 22
    22 = 126
    23 = 71
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
